# transfer_marbles_images.pptx - "Fix last image in transfer sample"
#
# 1) The auto date placeholder (datetimeFigureOut field) cached text moves
#    from 5/20/20 -> 5/21/20 on the slide master and every slide layout.
# 2) On the last slide (the "Org2" / post-transfer state diagram) the
#    marble's private-collection value changes from 110 -> 100 in the three
#    spots that still showed the old value:
#      - Channel World State "Can 5"  : hash(S:marble1) -> hash(100)
#      - Org1 private data collection : S:marble1 value  -> 100
#      - Channel World State "Can 37" : hash(S:marble1) -> hash(100)

$p = $ppt.ActivePresentation

function Set-DateFieldText($shapes, $newText) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $sh = $shapes.Item($i)
        $isDatePlaceholder = $false

        if ($sh.Type -eq 14) {
            try {
                if ($sh.PlaceholderFormat.Type -eq 16) {
                    $isDatePlaceholder = $true
                }
            } catch {
                $isDatePlaceholder = $false
            }
        }

        if ((-not $isDatePlaceholder) -and ($sh.Name -like "Date*")) {
            $isDatePlaceholder = $true
        }

        if ($isDatePlaceholder -and $sh.HasTextFrame) {
            $tr = $sh.TextFrame.TextRange
            if ($tr.Text -ne $newText) {
                $tr.Text = $newText
            }
        }
    }
}

# --- 1) Slide master + every slide layout ---------------------------------
$master = $p.SlideMaster
Set-DateFieldText $master.Shapes "5/21/20"

$layouts = $master.CustomLayouts
for ($li = 1; $li -le $layouts.Count; $li++) {
    $layout = $layouts.Item($li)
    Set-DateFieldText $layout.Shapes "5/21/20"
}

# --- 2) Last slide text fixes ----------------------------------------------
function Replace-RunText($shape, $oldText, $newText) {
    $tr = $shape.TextFrame.TextRange
    $fullText = $tr.Text
    $idx = $fullText.IndexOf($oldText)
    while ($idx -ge 0) {
        $sub = $tr.Characters($idx + 1, $oldText.Length)
        $sub.Text = $newText
        $fullText = $tr.Text
        $searchFrom = $idx + $newText.Length
        if ($searchFrom -ge $fullText.Length) {
            $idx = -1
        } else {
            $idx = $fullText.IndexOf($oldText, $searchFrom)
        }
    }
}

$lastSlide = $p.Slides.Item($p.Slides.Count)

$can5 = $lastSlide.Shapes.Item("Can 5")
Replace-RunText $can5 "hash(110)" "hash(100)"

$can20 = $lastSlide.Shapes.Item("Can 20")
Replace-RunText $can20 " 110" " 100"

$can37 = $lastSlide.Shapes.Item("Can 37")
Replace-RunText $can37 "hash(110)" "hash(100)"
